# Rename SKU prefix "ABQY" -> "NSQL" for all 100 product rows (A2:A101).
# SKU codes are zero-padded 5-digit sequence numbers: ABQY00001 .. ABQY00100
# become NSQL00001 .. NSQL00100, row-for-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 101; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $oldValue = [string]$cell.Value2
    $newValue = $oldValue -replace '^ABQY', 'NSQL'
    $cell.Value2 = $newValue
}
